$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric cell values (rows 2-9, columns D..AJ) per the IFRS data correction.
$ws.Range("AA2").Value = 54.37
$ws.Range("AB2").Value = 738.9
$ws.Range("AC2").Value = 2371
$ws.Range("AD2").Value = 13.79
$ws.Range("AE2").Value = 46591
$ws.Range("AF2").Value = 0.7
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 0.92
$ws.Range("AI2").Value = 12.49
$ws.Range("AJ2").Value = 25800000
$ws.Range("D2").Value = 12105
$ws.Range("E2").Value = -244
$ws.Range("F2").Value = -244
$ws.Range("G2").Value = 704
$ws.Range("H2").Value = 587
$ws.Range("I2").Value = 612
$ws.Range("J2").Value = -25
$ws.Range("K2").Value = 18460
$ws.Range("L2").Value = 6502
$ws.Range("M2").Value = 11958
$ws.Range("N2").Value = 11867
$ws.Range("O2").Value = 91
$ws.Range("P2").Value = 1290
$ws.Range("Q2").Value = 239
$ws.Range("R2").Value = 956
$ws.Range("S2").Value = -1093
$ws.Range("T2").Value = 761
$ws.Range("U2").Value = -522
$ws.Range("V2").Value = 4299
$ws.Range("W2").Value = -2.01
$ws.Range("X2").Value = 4.85
$ws.Range("Y2").Value = 5.29
$ws.Range("Z2").Value = 3.13
$ws.Range("AA3").Value = 44.75
$ws.Range("AB3").Value = 809.4
$ws.Range("AC3").Value = 3638
$ws.Range("AD3").Value = 9.789999999999999
$ws.Range("AE3").Value = 45684
$ws.Range("AF3").Value = 0.78
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.4
$ws.Range("AI3").Value = 13.57
$ws.Range("AJ3").Value = 25800000
$ws.Range("D3").Value = 11619
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = 1127
$ws.Range("H3").Value = 900
$ws.Range("I3").Value = 939
$ws.Range("J3").Value = -39
$ws.Range("K3").Value = 16842
$ws.Range("L3").Value = 5206
$ws.Range("M3").Value = 11636
$ws.Range("N3").Value = 11636
$ws.Range("P3").Value = 1290
$ws.Range("Q3").Value = 995
$ws.Range("R3").Value = 841
$ws.Range("S3").Value = -614
$ws.Range("T3").Value = 841
$ws.Range("U3").Value = 154
$ws.Range("V3").Value = 3171
$ws.Range("W3").Value = 0.22
$ws.Range("X3").Value = 7.74
$ws.Range("Y3").Value = 7.99
$ws.Range("Z3").Value = 5.1
$ws.Range("AA4").Value = 37.72
$ws.Range("AB4").Value = 769.14
$ws.Range("AC4").Value = -1688
$ws.Range("AD4").Value = -18.16
$ws.Range("AE4").Value = 43646
$ws.Range("AF4").Value = 0.7
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 0.98
$ws.Range("AI4").Value = -17.55
$ws.Range("AJ4").Value = 25800000
$ws.Range("D4").Value = 11107
$ws.Range("E4").Value = 298
$ws.Range("F4").Value = 298
$ws.Range("G4").Value = -504
$ws.Range("H4").Value = -435
$ws.Range("I4").Value = -435
$ws.Range("K4").Value = 15309
$ws.Range("L4").Value = 4193
$ws.Range("M4").Value = 11117
$ws.Range("N4").Value = 11117
$ws.Range("P4").Value = 1290
$ws.Range("Q4").Value = 655
$ws.Range("R4").Value = -459
$ws.Range("S4").Value = -287
$ws.Range("T4").Value = 312
$ws.Range("U4").Value = 343
$ws.Range("V4").Value = 3014
$ws.Range("W4").Value = 2.68
$ws.Range("X4").Value = -3.92
$ws.Range("Y4").Value = -3.83
$ws.Range("Z4").Value = -2.71
$ws.Range("AA5").Value = 37.61
$ws.Range("AB5").Value = 832.05
$ws.Range("AC5").Value = 3459
$ws.Range("AD5").Value = 18.77
$ws.Range("AE5").Value = 46828
$ws.Range("AF5").Value = 1.39
$ws.Range("AG5").Value = 800
$ws.Range("AH5").Value = 1.23
$ws.Range("AI5").Value = 22.84
$ws.Range("AJ5").Value = 25800000
$ws.Range("D5").Value = 11595
$ws.Range("E5").Value = 1111
$ws.Range("F5").Value = 1111
$ws.Range("G5").Value = 1107
$ws.Range("H5").Value = 892
$ws.Range("I5").Value = 892
$ws.Range("K5").Value = 16413
$ws.Range("L5").Value = 4486
$ws.Range("M5").Value = 11927
$ws.Range("N5").Value = 11927
$ws.Range("P5").Value = 1290
$ws.Range("Q5").Value = 1888
$ws.Range("R5").Value = -2389
$ws.Range("S5").Value = -282
$ws.Range("T5").Value = 629
$ws.Range("U5").Value = 1260
$ws.Range("V5").Value = 2810
$ws.Range("W5").Value = 9.58
$ws.Range("X5").Value = 7.7
$ws.Range("Y5").Value = 7.74
$ws.Range("Z5").Value = 5.63
$ws.Range("AA6").Value = 23.38
$ws.Range("AB6").Value = 989.97
$ws.Range("AC6").Value = 8330
$ws.Range("AD6").Value = 4.93
$ws.Range("AE6").Value = 54099
$ws.Range("AF6").Value = 0.76
$ws.Range("AG6").Value = 1700
$ws.Range("AH6").Value = 4.14
$ws.Range("AI6").Value = 20.15
$ws.Range("AJ6").Value = 25800000
$ws.Range("D6").Value = 13717
$ws.Range("E6").Value = 2107
$ws.Range("F6").Value = 2107
$ws.Range("G6").Value = 2752
$ws.Range("H6").Value = 2149
$ws.Range("I6").Value = 2149
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 3221
$ws.Range("M6").Value = 13779
$ws.Range("N6").Value = 13779
$ws.Range("P6").Value = 1290
$ws.Range("Q6").Value = 2732
$ws.Range("R6").Value = -1777
$ws.Range("S6").Value = -2204
$ws.Range("T6").Value = 965
$ws.Range("U6").Value = 1767
$ws.Range("V6").Value = 812
$ws.Range("W6").Value = 15.36
$ws.Range("X6").Value = 15.67
$ws.Range("Y6").Value = 16.72
$ws.Range("Z6").Value = 12.86
$ws.Range("AA7").Value = 21.37
$ws.Range("AC7").Value = 6587
$ws.Range("AD7").Value = 5.93
$ws.Range("AE7").Value = 58942
$ws.Range("AF7").Value = 0.66
$ws.Range("AG7").Value = 1700
$ws.Range("AH7").Value = 4.35
$ws.Range("AI7").Value = 25.81
$ws.Range("D7").Value = 13066
$ws.Range("E7").Value = 1854
$ws.Range("G7").Value = 2245
$ws.Range("H7").Value = 1703
$ws.Range("I7").Value = 1700
$ws.Range("K7").Value = 18269
$ws.Range("L7").Value = 3216
$ws.Range("M7").Value = 15053
$ws.Range("N7").Value = 15012
$ws.Range("P7").Value = 1290
$ws.Range("Q7").Value = 2416
$ws.Range("R7").Value = -1076
$ws.Range("S7").Value = -417
$ws.Range("T7").Value = 1098
$ws.Range("U7").Value = 1217
$ws.Range("W7").Value = 14.19
$ws.Range("X7").Value = 13.03
$ws.Range("Y7").Value = 11.81
$ws.Range("Z7").Value = 9.66
$ws.Range("AA8").Value = 18.81
$ws.Range("AC8").Value = 7192
$ws.Range("AD8").Value = 5.43
$ws.Range("AE8").Value = 64468
$ws.Range("AF8").Value = 0.61
$ws.Range("AG8").Value = 1960
$ws.Range("AH8").Value = 5.02
$ws.Range("AI8").Value = 27.25
$ws.Range("D8").Value = 13740
$ws.Range("E8").Value = 2022
$ws.Range("G8").Value = 2420
$ws.Range("H8").Value = 1856
$ws.Range("I8").Value = 1856
$ws.Range("K8").Value = 19534
$ws.Range("L8").Value = 3093
$ws.Range("M8").Value = 16442
$ws.Range("N8").Value = 16420
$ws.Range("P8").Value = 1290
$ws.Range("Q8").Value = 2291
$ws.Range("R8").Value = -1343
$ws.Range("S8").Value = -626
$ws.Range("T8").Value = 1136
$ws.Range("U8").Value = 971
$ws.Range("W8").Value = 14.72
$ws.Range("X8").Value = 13.51
$ws.Range("Y8").Value = 11.81
$ws.Range("Z8").Value = 9.82
$ws.Range("AA9").Value = 18.77
$ws.Range("AC9").Value = 8255
$ws.Range("AD9").Value = 4.73
$ws.Range("AE9").Value = 70907
$ws.Range("AF9").Value = 0.55
$ws.Range("AG9").Value = 1960
$ws.Range("AH9").Value = 5.02
$ws.Range("AI9").Value = 23.74
$ws.Range("D9").Value = 14439
$ws.Range("E9").Value = 2363
$ws.Range("G9").Value = 2768
$ws.Range("H9").Value = 2129
$ws.Range("I9").Value = 2130
$ws.Range("K9").Value = 21463
$ws.Range("L9").Value = 3392
$ws.Range("M9").Value = 18071
$ws.Range("N9").Value = 18060
$ws.Range("P9").Value = 1290
$ws.Range("Q9").Value = 2576
$ws.Range("R9").Value = -1510
$ws.Range("S9").Value = -280
$ws.Range("T9").Value = 1164
$ws.Range("U9").Value = 1020
$ws.Range("W9").Value = 16.37
$ws.Range("X9").Value = 16.37
$ws.Range("Y9").Value = 12.35
$ws.Range("Z9").Value = 10.39

# Remove cells that no longer exist in the corrected data (fully cleared, not just blanked).
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
